{"js": "// Apply the title/description/bullet rewrites described by the diff.\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    from: \"Play Fu Fortunes Megaways for Free | Game Review\",\n    to: \"Play Fu Fortunes Megaways Free - Exciting Asian Themed Slot\"\n  },\n  {\n    from: \"Megaways mechanism allows up to 46,656 ways to win\",\n    to: \"Megaways mechanism with up to 46,656 ways to win\"\n  },\n  {\n    from: \"Bonus functions increase chances of winning\",\n    to: \"Medium to high volatility for exciting gameplay\"\n  },\n  {\n    from: \"Four jackpots offer larger payouts\",\n    to: \"Includes wild symbol, jackpots, and free spins\"\n  },\n  {\n    from: \"Free spins offer multipliers up to 20x\",\n    to: \"High payout potential of up to 9710 times your bet\"\n  },\n  {\n    from: \"Medium to high volatility may deter some players\",\n    to: \"Limited bet range with minimum bet of 20 cents\"\n  },\n  {\n    from: \"Maximum bet limit of 20 euros may not appeal to high rollers\",\n    to: \"Payout percentage slightly below average for online slots\"\n  },\n  {\n    from: \"Read our review of Fu Fortunes Megaways, an online slot with 6 reels and up to 46,656 ways to win. Play for free and discover bonus functions and jackpots.\",\n    to: \"Discover Fu Fortunes Megaways, a thrilling Asian themed slot game with free spins and big winnings. Play for free!\"\n  }\n];\n\nfor (const { from, to } of replacements) {\n  const results = body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the title/description/bullet rewrites described by the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ From = \"Play Fu Fortunes Megaways for Free | Game Review\"; To = \"Play Fu Fortunes Megaways Free - Exciting Asian Themed Slot\" },\n    @{ From = \"Megaways mechanism allows up to 46,656 ways to win\"; To = \"Megaways mechanism with up to 46,656 ways to win\" },\n    @{ From = \"Bonus functions increase chances of winning\"; To = \"Medium to high volatility for exciting gameplay\" },\n    @{ From = \"Four jackpots offer larger payouts\"; To = \"Includes wild symbol, jackpots, and free spins\" },\n    @{ From = \"Free spins offer multipliers up to 20x\"; To = \"High payout potential of up to 9710 times your bet\" },\n    @{ From = \"Medium to high volatility may deter some players\"; To = \"Limited bet range with minimum bet of 20 cents\" },\n    @{ From = \"Maximum bet limit of 20 euros may not appeal to high rollers\"; To = \"Payout percentage slightly below average for online slots\" },\n    @{ From = \"Read our review of Fu Fortunes Megaways, an online slot with 6 reels and up to 46,656 ways to win. Play for free and discover bonus functions and jackpots.\"; To = \"Discover Fu Fortunes Megaways, a thrilling Asian themed slot game with free spins and big winnings. Play for free!\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.From\n    $find.Replacement.Text = $r.To\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
